$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1,1).Value = 1.732704007046913
$ws.Cells.Item(1,2).Value = 0.1945306715051764
$ws.Cells.Item(1,3).Value = 0.672838158291254
$ws.Cells.Item(1,4).Value = 0.7034277224914169
$ws.Cells.Item(1,5).Value = 1.570796292848413
$ws.Cells.Item(1,6).Value = 0.161907685808031
$ws.Cells.Item(2,1).Value = 1.778290465355663
$ws.Cells.Item(2,2).Value = 0.1934420031981677
$ws.Cells.Item(2,3).Value = 0.6758680428538398
$ws.Cells.Item(2,4).Value = 0.7014865303281086
$ws.Cells.Item(2,5).Value = 1.570796289207822
$ws.Cells.Item(2,6).Value = 0.2074941448287797
$ws.Cells.Item(3,1).Value = 1.982548021645145
$ws.Cells.Item(3,2).Value = 0.1885640474137254
$ws.Cells.Item(3,3).Value = 0.6894439344574702
$ws.Cells.Item(3,4).Value = 0.6927887024582571
$ws.Cells.Item(3,5).Value = 1.570796272895559
$ws.Cells.Item(3,6).Value = 0.411751704308488
$ws.Cells.Item(4,1).Value = 2.269919010210249
$ws.Cells.Item(4,2).Value = 0.1817012267601512
$ws.Cells.Item(4,3).Value = 0.708543924987362
$ws.Cells.Item(4,4).Value = 0.680551684456407
$ws.Cells.Item(4,5).Value = 1.570796249945755
$ws.Cells.Item(4,6).Value = 0.6991226973619369
$ws.Cells.Item(5,1).Value = 2.474176566499732
$ws.Cells.Item(5,2).Value = 0.1768232709757089
$ws.Cells.Item(5,3).Value = 0.7221198165909923
$ws.Cells.Item(5,4).Value = 0.6718538565865555
$ws.Cells.Item(5,5).Value = 1.570796233633492
$ws.Cells.Item(5,6).Value = 0.9033802568416449
$ws.Cells.Item(6,1).Value = 2.519763024808483
$ws.Cells.Item(6,2).Value = 0.1757346026687003
$ws.Cells.Item(6,3).Value = 0.7251497011535781
$ws.Cells.Item(6,4).Value = 0.6699126644232472
$ws.Cells.Item(6,5).Value = 1.570796229992901
$ws.Cells.Item(6,6).Value = 0.9489667158623929
$ws.Cells.Item(7,1).Value = 1.900875918751935
$ws.Cells.Item(7,2).Value = 0.05335115204042857
$ws.Cells.Item(7,3).Value = 1.199938781730575
$ws.Cells.Item(7,4).Value = 0.3175060894984869
$ws.Cells.Item(7,5).Value = 1.570796342143612
$ws.Cells.Item(7,6).Value = 0.3300795436886814
$ws.Cells.Item(8,1).Value = 1.639788217465948
$ws.Cells.Item(8,2).Value = 0.07454629022105426
$ws.Cells.Item(8,3).Value = 1.120803451180303
$ws.Cells.Item(8,4).Value = 0.3754463736539357
$ws.Cells.Item(8,5).Value = 1.570796332944535
$ws.Cells.Item(8,6).Value = 0.06899185004407737
$ws.Cells.Item(9,1).Value = 1.361539650321506
$ws.Cells.Item(9,2).Value = 0.09713455002854179
$ws.Cells.Item(9,3).Value = 1.036466685305622
$ws.Cells.Item(9,4).Value = 0.437194977539862
$ws.Cells.Item(9,5).Value = 1.570796323140818
$ws.Cells.Item(9,6).Value = -0.2092567089567264
$ws.Cells.Item(10,1).Value = 1.083291083177064
$ws.Cells.Item(10,2).Value = 0.1197228098360293
$ws.Cells.Item(10,3).Value = 0.9521299194309416
$ws.Cells.Item(10,4).Value = 0.4989435814257883
$ws.Cells.Item(10,5).Value = 1.570796313337101
$ws.Cells.Item(10,6).Value = -0.4875052679575304
$ws.Cells.Item(11,1).Value = 0.8222033818910772
$ws.Cells.Item(11,2).Value = 0.140917948016655
$ws.Cells.Item(11,3).Value = 0.8729945888806701
$ws.Cells.Item(11,4).Value = 0.5568838655812371
$ws.Cells.Item(11,5).Value = 1.570796304138024
$ws.Cells.Item(11,6).Value = -0.748592961602134
$ws.Cells.Item(12,1).Value = 0.5937899691995894
$ws.Cells.Item(12,2).Value = 0.1594605826197359
$ws.Cells.Item(12,3).Value = 0.8037627911880737
$ws.Cells.Item(12,4).Value = 0.6075731089698566
$ws.Cells.Item(12,5).Value = 1.570796296090182
$ws.Cells.Item(12,6).Value = -0.977006367608534
$ws.Cells.Item(13,1).Value = 0.410269381593821
$ws.Cells.Item(13,2).Value = 0.1743588110469463
$ws.Cells.Item(13,3).Value = 0.7481379483041316
$ws.Cells.Item(13,4).Value = 0.6482997879435473
$ws.Cells.Item(13,5).Value = 1.570796289624078
$ws.Cells.Item(13,6).Value = -1.160526949843115
$ws.Cells.Item(14,1).Value = 0.2789178261977563
$ws.Cells.Item(14,2).Value = 0.1850219497284968
$ws.Cells.Item(14,3).Value = 0.7083254688063935
$ws.Cells.Item(14,4).Value = 0.6774491749365862
$ws.Cells.Item(14,5).Value = 1.570796284996081
$ws.Cells.Item(14,6).Value = -1.291878501394849
$ws.Cells.Item(15,1).Value = 0.2004217376457329
$ws.Cells.Item(15,2).Value = 0.1913942737993129
$ws.Cells.Item(15,3).Value = 0.6845334101078362
$ws.Cells.Item(15,4).Value = 0.6948689371597545
$ws.Cells.Item(15,5).Value = 1.570796282230376
$ws.Cells.Item(15,6).Value = -1.370374587649487
$ws.Cells.Item(16,1).Value = 0.1672303349600339
$ws.Cells.Item(16,2).Value = 0.194088756775214
$ws.Cells.Item(16,3).Value = 0.6744731406657185
$ws.Cells.Item(16,4).Value = 0.7022347352944619
$ws.Cells.Item(16,5).Value = 1.570796281060922
$ws.Cells.Item(16,6).Value = -1.403565989363758
$ws.Cells.Item(17,1).Value = 0.1619081784284635
$ws.Cells.Item(17,2).Value = 0.1945208102290914
$ws.Cells.Item(17,3).Value = 0.6728600021904421
$ws.Cells.Item(17,4).Value = 0.7034158221868742
$ws.Cells.Item(17,5).Value = 1.570796280873402
$ws.Cells.Item(17,6).Value = -1.40888814573956
